$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 39; this shifts the existing rows 39..81
# down to 40..82, preserving all of their data (and formatting) intact.
$ws.Rows.Item(39).Insert()

# Populate the newly-inserted row 39 with the new weekly record.
$ws.Cells.Item(39, 1).Value = 10
$ws.Cells.Item(39, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(39, 3).Value = "La Araucanía"
$ws.Cells.Item(39, 4).Value = 45033
$ws.Cells.Item(39, 5).Value = 9
$ws.Cells.Item(39, 6).Value = 100112010
$ws.Cells.Item(39, 7).Value = "Achicoria"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 55
$ws.Cells.Item(39, 11).Value = 10000
$ws.Cells.Item(39, 12).Value = 10000
$ws.Cells.Item(39, 13).Value = 10000
$ws.Cells.Item(39, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(39, 15).Value = "Región Metropolitana"
$ws.Cells.Item(39, 16).Value = 556
$ws.Cells.Item(39, 17).Value = 18
$ws.Cells.Item(39, 18).Value = "Hortaliza"
